$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new prompt-library rows are appended in column F (rows 10 and 11),
# following the same "title/description" shorthand used by every other
# entry in the sheet. Both are long, so they get word-wrap plus a taller
# row to show the full text (mirrors how the existing F1:F9 entries look
# once Excel wraps long prompts).
$verbText = "动词/提供20个不同的常用动词，in the following format:     - Keyword 1     - Keyword 2     - Keyword 3`n"
$synonymText = "日语近义词/提供20个和主题内容相似的日语单词，提供例句和中文翻译，讲解语法，具体说明使用上的差别。in the following format:     - 段落 1     - 段落 2     - 段落 3`n"

$cellF10 = $ws.Range("F10")
$cellF10.Value = $verbText
$cellF10.VerticalAlignment = -4108 # xlVAlignCenter
$cellF10.WrapText = $true

$cellF11 = $ws.Range("F11")
$cellF11.Value = $synonymText
$cellF11.VerticalAlignment = -4108 # xlVAlignCenter
$cellF11.WrapText = $true

# Grow the rows enough to show the wrapped text, matching how Excel sized
# them once the long prompt text was wrapped in that column.
$ws.Rows.Item(10).RowHeight = 177
$ws.Rows.Item(11).RowHeight = 231.75

# Leave the selection/scroll position where the author ended up: looking
# at the newly-added rows near the bottom of the sheet, with F12 (the
# empty cell right after the new content) active.
$ws.Range("F12").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
